# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    populated with the per-fund holding detail for the 2022-Q1 quarter
#    (same column layout as the other per-quarter sheets: 基金代码, 基金名称,
#    基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名).
# 2. Add a matching summary row to the top of the "总计" sheet's data table
#    (date/count/value), shifting the existing rows down and renumbering the
#    leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B:G hold text (fund codes must keep leading zeros, and the
# percentages/amounts are stored as plain text in this workbook's format).
$q1.Range("B2:G3").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "005947"
$q1.Range("C2").Value = "德邦民裕进取量化精选灵活配置混合A"
$q1.Range("D2").Value = "0.53"
$q1.Range("E2").Value = "94.44"
$q1.Range("F2").Value = "6.75"
$q1.Range("G2").Value = "0.0358"
$q1.Range("H2").Value = 6

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "005948"
$q1.Range("C3").Value = "德邦民裕进取量化精选灵活配置混合C"
$q1.Range("D3").Value = "0.09"
$q1.Range("E3").Value = "94.44"
$q1.Range("F3").Value = "6.75"
$q1.Range("G3").Value = "0.0061"
$q1.Range("H3").Value = 6

# Drop the "@" number-format footprint from the data cells (keep them stored
# as Text, but back to the sheet's default/unstyled look) by pasting formats
# from a never-touched cell.
$q1.Range("Z100").Copy()
$q1.Range("B2:G3").PasteSpecial(-4122)

# Style the header row (bold + border) the same way every other per-quarter
# sheet's header row is styled.
$otherHeader = $wb.Worksheets.Item("2021-Q4").Range("B1")
$otherHeader.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Insert the new quarter's row into "总计", above the existing data
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Rows.Item(2).ClearFormats()

# Match the index column's look (style) used by the rest of the table.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

# Renumber the leading index column for the rows pushed down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
